# MonitoramentoCargaOBM.xlsx update
# - Update the April (row 5) totals with revised figures
# - Add the May (row 6) breakdown figures
# - Add reviewer comments on A5 and A6 documenting the monitoring dates
# - Leave the cursor selection where the author left it (C14)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 ("ABRIL 23") revised figures -----------------------------------
$ws.Range("C5").Value = 1155
$ws.Range("D5").Value = 93
$ws.Range("E5").Value = 93
$ws.Range("F5").Value = 545
$ws.Range("G5").Value = 650
$ws.Range("H5").Value = 596

# --- Row 6 ("MAIO") new figures --------------------------------------------
$ws.Range("C6").Value = 1385
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 79
$ws.Range("F6").Value = 305
$ws.Range("G6").Value = 850
$ws.Range("H6").Value = 1223

# --- Reviewer comments ------------------------------------------------------
$excel.UserName = "ips.nt0001@gmail.com"

$ws.Range("A5").AddComment("`nIniciado 17/04")
$ws.Range("A6").AddComment("`nAté dia 21/05")

# --- Restore the selection left by the author ------------------------------
$ws.Range("C14").Select()
